$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3493
$ws.Range("F5").Value = 8176
$ws.Range("F8").Value = 2143
$ws.Range("F10").Value = 168
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 31
$ws.Range("F15").Value = 751
$ws.Range("F17").Value = 552
$ws.Range("F20").Value = 1376
$ws.Range("F21").Value = 6932
$ws.Range("F22").Value = 132
$ws.Range("F23").Value = 54328
$ws.Range("F24").Value = 4262
$ws.Range("F25").Value = 9
$ws.Range("F27").Value = 1024
$ws.Range("F28").Value = 818
$ws.Range("F29").Value = 393
$ws.Range("F30").Value = 78
$ws.Range("F31").Value = 860
$ws.Range("F34").Value = 2053
$ws.Range("F36").Value = 573
$ws.Range("F38").Value = 846
$ws.Range("F39").Value = 1116
$ws.Range("F40").Value = 565
$ws.Range("F42").Value = 176
$ws.Range("F43").Value = 1055
$ws.Range("F44").Value = 691
$ws.Range("F45").Value = 130
$ws.Range("F47").Value = 123
$ws.Range("F48").Value = 33

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 26
$ws.Range("F14").Value = 33
$ws.Range("F15").Value = 168
$ws.Range("F16").Value = 7407
$ws.Range("G16").Value = 680
$ws.Range("F17").Value = 99
$ws.Range("F29").Value = 76

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2257
$ws.Range("F9").Value = 9297
$ws.Range("F10").Value = 1602
$ws.Range("F11").Value = 156
$ws.Range("F12").Value = 67
$ws.Range("F15").Value = 124

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2257
$ws.Range("F4").Value = 8176
$ws.Range("F6").Value = 1602
$ws.Range("F7").Value = 156
$ws.Range("F9").Value = 2143
$ws.Range("F13").Value = 168
$ws.Range("F15").Value = 7
$ws.Range("F16").Value = 552
$ws.Range("F18").Value = 6932
$ws.Range("F19").Value = 132
$ws.Range("F20").Value = 54329
$ws.Range("F23").Value = 26
$ws.Range("F25").Value = 4262
$ws.Range("F26").Value = 9
$ws.Range("F28").Value = 818
$ws.Range("F29").Value = 393
$ws.Range("F30").Value = 78
$ws.Range("F32").Value = 2053
$ws.Range("F33").Value = 33
$ws.Range("F34").Value = 846
$ws.Range("F35").Value = 1117
$ws.Range("F38").Value = 1055
$ws.Range("F40").Value = 691
$ws.Range("F42").Value = 130
$ws.Range("F44").Value = 123
$ws.Range("F45").Value = 76
$ws.Range("F46").Value = 33
